$wb = $excel.ActiveWorkbook

# --- Phase2 sheet: update probabilities and selection ---
$ws2 = $wb.Worksheets.Item("Phase2")
$ws2.Activate()

$ws2.Range("C2").Value = 0.1
$ws2.Range("C3").Value = 0.1
$ws2.Range("C4").Value = 0.05
$ws2.Range("C5").Value = 0.05

$ws2.Range("C2:C5").Select() | Out-Null

# --- Phase1 sheet: update probabilities and selection ---
$ws1 = $wb.Worksheets.Item("Phase1")
$ws1.Activate()

$ws1.Range("B2").Value = 0.1
$ws1.Range("B3").Value = 0.1
$ws1.Range("B4").Value = 0.05
$ws1.Range("B5").Value = 0.05

$ws1.Range("B2:B5").Select() | Out-Null
